$d = $word.ActiveDocument

# Remove the standalone "in progress" paragraph (bold/italic, FirstParagraph
# style) that directly follows the "F2025" date paragraph. Deleting the
# paragraph's Range (rather than just its text) also removes the paragraph
# mark, so the document collapses back to the state before that paragraph
# was added.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "in progress") {
        $p.Range.Delete()
        break
    }
}
